# Fixed error with initial software defaulting to 1 without input
#
# Appends two new log rows to "repair_notes" (rows 12-13) and three new
# rows to "new_inventory" (rows 11-13) describing the same check-ins for
# host LAPTOP-K9FQV11C. The "initial software" flags (ArcMap/FoxIT
# Pro/BlueBeam/LanSweeper columns Q/U/X/Y) now default to "0" instead of
# "1" when nothing was supplied, except where a real answer ("1") was
# captured for row 13's ArcMap column.

function Set-TextCell {
    param($ws, $row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    # Start from the plain/unstyled look (matches the other un-styled data
    # rows already in this workbook, e.g. repair_notes rows 6-11) so that
    # forcing a text number-format below doesn't also drag in inherited
    # column formatting (e.g. the wrapText alignment on repair_notes'
    # column style). Force text storage (matches the rest of the sheet,
    # which stores dates/times/flags as plain strings rather than
    # numbers/dates), then drop back to the default look once more.
    $cell.Style = "Normal"
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-BlankCell {
    param($ws, $row, $col)
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = ""
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# repair_notes: two new check-in rows for LAPTOP-K9FQV11C
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Item("repair_notes")

Set-TextCell $notes 12 1  "LAPTOP-K9FQV11C"
Set-TextCell $notes 12 2  "Windows-10-10.0.19041-SP0"
Set-TextCell $notes 12 3  "Intel(R) Core(TM) i5-7200U CPU @ 2.50GHz"
Set-TextCell $notes 12 4  "16 GB"
Set-TextCell $notes 12 5  "192.168.1.130"
Set-TextCell $notes 12 6  "12-20-2020"
Set-TextCell $notes 12 7  "08:39"
Set-TextCell $notes 12 8  "home"
Set-TextCell $notes 12 9  "This is my personal laptop"
Set-TextCell $notes 12 10 "Dan"

Set-TextCell $notes 13 1  "LAPTOP-K9FQV11C"
Set-TextCell $notes 13 2  "Windows-10-10.0.19041-SP0"
Set-TextCell $notes 13 3  "Intel(R) Core(TM) i5-7200U CPU @ 2.50GHz"
Set-TextCell $notes 13 4  "16 GB"
Set-TextCell $notes 13 5  "192.168.1.130"
Set-TextCell $notes 13 6  "12-20-2020"
Set-TextCell $notes 13 7  "09:08"
Set-TextCell $notes 13 8  "PY_VAR1"
Set-TextCell $notes 13 10 "PY_VAR0"

# ---------------------------------------------------------------------
# new_inventory: matching inventory rows for LAPTOP-K9FQV11C.
# Software flags (Q/U/X/Y = ArcMap/FoxIT Pro/BlueBeam/LanSweeper) default
# to "0" (previously incorrectly defaulted to "1" with no input).
# ---------------------------------------------------------------------
$inv = $wb.Worksheets.Item("new_inventory")

Set-TextCell $inv 11 1  "LAPTOP-K9FQV11C"
Set-TextCell $inv 11 2  "Y"
Set-TextCell $inv 11 8  "PY_VAR0"
Set-TextCell $inv 11 9  "PY_VAR1"
Set-TextCell $inv 11 17 "0"
Set-TextCell $inv 11 21 "0"
Set-TextCell $inv 11 24 "0"
Set-TextCell $inv 11 25 "0"

Set-TextCell $inv 12 1  "LAPTOP-K9FQV11C"
Set-TextCell $inv 12 2  "Y"
Set-TextCell $inv 12 8  "PY_VAR0"
Set-TextCell $inv 12 9  "PY_VAR1"
Set-TextCell $inv 12 17 "0"
Set-TextCell $inv 12 21 "0"
Set-TextCell $inv 12 24 "0"
Set-TextCell $inv 12 25 "0"

Set-TextCell $inv 13 1  "LAPTOP-K9FQV11C"
Set-TextCell $inv 13 2  "Y"
Set-BlankCell $inv 13 3
Set-BlankCell $inv 13 4
Set-BlankCell $inv 13 5
Set-BlankCell $inv 13 6
Set-TextCell $inv 13 8  "PY_VAR0"
Set-TextCell $inv 13 9  "PY_VAR1"
Set-TextCell $inv 13 17 "1"
Set-TextCell $inv 13 21 "0"
Set-TextCell $inv 13 24 "0"
Set-TextCell $inv 13 25 "0"
